$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New activity entries for 08/05/2019 (date serial 43593)
$date = Get-Date -Year 2019 -Month 5 -Day 8 -Hour 0 -Minute 0 -Second 0

$ws.Range("A75").Value = "Giovanni"
$ws.Range("B75").Value = "Documenti di progetto"
$ws.Range("C75").Value = 45
$ws.Range("D75").Value = $date

$ws.Range("A76").Value = "Hristina"
$ws.Range("B76").Value = "Documenti di progetto"
$ws.Range("C76").Value = 45
$ws.Range("D76").Value = $date

$ws.Range("A77").Value = "Luca"
$ws.Range("B77").Value = "Documenti di progetto"
$ws.Range("C77").Value = 45
$ws.Range("D77").Value = $date

$ws.Range("A78").Value = "Viktorija"
$ws.Range("B78").Value = "Documenti di progetto"
$ws.Range("C78").Value = 45
$ws.Range("D78").Value = $date

$ws.Range("A79").Value = "Hristina"
$ws.Range("B79").Value = "Sviluppo"
$ws.Range("C79").Value = 60
$ws.Range("D79").Value = $date

$ws.Range("A80").Value = "Luca"
$ws.Range("B80").Value = "Sviluppo"
$ws.Range("C80").Value = 60
$ws.Range("D80").Value = $date

$ws.Range("A81").Value = "Viktorija"
$ws.Range("B81").Value = "Sviluppo"
$ws.Range("C81").Value = 60
$ws.Range("D81").Value = $date

# Update the view state to match where the user navigated/selected
$excel.ActiveWindow.ScrollRow = 55
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("D81").Select()

$wb.Save()
